# Fix issue security department field name in template, and join
# a couple of runs that were previously split mid-word back into
# single runs (the wording itself is unchanged, only the run
# splitting and the template field name are fixed).

$d = $word.ActiveDocument

# 1) "Отсутстви" + "е" + " судимостей ..." -> single run with full text.
$d.Content.Find.Execute(
    "Отсутствие судимостей в отношении физических лиц (генеральный директор, участники юридического лица (c наибольшей долей участия, Поручитель)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Отсутствие судимостей в отношении физических лиц (генеральный директор, участники юридического лица (c наибольшей долей участия, Поручитель)",
    2) | Out-Null

# 2) Rename the merge-field: court_acts_info -> court_cases_info.
$d.Content.Find.Execute(
    "{issue.issuer_shareholders_participants_or_self_court_acts_info}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{issue.issuer_shareholders_participants_or_self_court_cases_info}",
    2) | Out-Null

# 3) "Информаци" + "я" + " о судебных разбирательствах ..." -> single run with full text.
$d.Content.Find.Execute(
    "Информация о судебных разбирательствах Принципала (в качестве ответчика), о находящихся в суде делах и принятых по ним судебным актам",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Информация о судебных разбирательствах Принципала (в качестве ответчика), о находящихся в суде делах и принятых по ним судебным актам",
    2) | Out-Null
